# Update countries & provincias Spain
# Refresh the COVID country stats table ("Pais" sheet) to the 19-Mar-2020
# 22:45 snapshot: new case numbers for several countries plus a handful of
# countries that changed rank (their row now shows a different country,
# with that country's updated figures) and the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value2 = "Datos actualizados a 19 de Marzo de 2020 a las 22:45"

# Row 9
$ws.Range("B9").Value2 = 13475
$ws.Range("C9").Value2 = 4216
$ws.Range("E9").Value2 = 13180
$ws.Range("G9").Value2 = 37
$ws.Range("H9").Value2 = 187

# Row 10
$ws.Range("D10").Value2 = 1295
$ws.Range("E10").Value2 = 9328
$ws.Range("F10").Value2 = 1122

# Row 12
$ws.Range("B12").Value2 = 4222
$ws.Range("C12").Value2 = 1107
$ws.Range("E12").Value2 = 4164

# Row 47
$ws.Range("A47").Value2 = "Ecuador"
$ws.Range("B47").Value2 = 260
$ws.Range("C47").Value2 = 92
$ws.Range("D47").Value2 = 1
$ws.Range("E47").Value2 = 256
$ws.Range("F47").Value2 = 2
$ws.Range("G47").Value2 = 0
$ws.Range("H47").Value2 = 3

# Row 48
$ws.Range("A48").Value2 = "Egipto"
$ws.Range("B48").Value2 = 256
$ws.Range("C48").Value2 = 46
$ws.Range("D48").Value2 = 42
$ws.Range("E48").Value2 = 207
$ws.Range("F48").Value2 = 0
$ws.Range("G48").Value2 = 1
$ws.Range("H48").Value2 = 7

# Row 49
$ws.Range("A49").Value2 = "Peru"
$ws.Range("B49").Value2 = 234
$ws.Range("C49").Value2 = 89
$ws.Range("D49").Value2 = 1
$ws.Range("E49").Value2 = 233
$ws.Range("F49").Value2 = 7
$ws.Range("H49").Value2 = 0

# Row 50
$ws.Range("A50").Value2 = "Filipinas"
$ws.Range("B50").Value2 = 217
$ws.Range("D50").Value2 = 8
$ws.Range("E50").Value2 = 192
$ws.Range("F50").Value2 = 1
$ws.Range("H50").Value2 = 17

# Row 51
$ws.Range("A51").Value2 = "Hong Kong"
$ws.Range("B51").Value2 = 208
$ws.Range("C51").Value2 = 15
$ws.Range("D51").Value2 = 98
$ws.Range("E51").Value2 = 106
$ws.Range("F51").Value2 = 4
$ws.Range("H51").Value2 = 4

# Row 58
$ws.Range("B58").Value2 = 144
$ws.Range("C58").Value2 = 4
$ws.Range("E58").Value2 = 126
$ws.Range("F58").Value2 = 12

# Row 111
$ws.Range("A111").Value2 = "Bolivia"
$ws.Range("C111").Value2 = 3

# Row 112
$ws.Range("A112").Value2 = "Guayana Francesa"
$ws.Range("C112").Value2 = 0

# Row 115
$ws.Range("A115").Value2 = "Montenegro"
$ws.Range("C115").Value2 = 5

# Row 117
$ws.Range("A117").Value2 = "Maldivas"
$ws.Range("B117").Value2 = 13
$ws.Range("C117").Value2 = 0
$ws.Range("E117").Value2 = 13

# Row 119
$ws.Range("A119").Value2 = "Honduras"
$ws.Range("C119").Value2 = 3
$ws.Range("D119").Value2 = 0
$ws.Range("E119").Value2 = 12

# Row 120
$ws.Range("A120").Value2 = "Nigeria"
$ws.Range("B120").Value2 = 12
$ws.Range("C120").Value2 = 4
$ws.Range("D120").Value2 = 1
$ws.Range("F120").Value2 = 0

# Row 121
$ws.Range("A121").Value2 = "Paraguay"
$ws.Range("C121").Value2 = 0
$ws.Range("F121").Value2 = 1

# Row 122
$ws.Range("A122").Value2 = "Ghana"
$ws.Range("C122").Value2 = 4

# Row 123
$ws.Range("A123").Value2 = "Ruanda"
$ws.Range("C123").Value2 = 0
$ws.Range("E123").Value2 = 11
$ws.Range("H123").Value2 = 0

# Row 124
$ws.Range("A124").Value2 = "Cuba"
$ws.Range("B124").Value2 = 11
$ws.Range("C124").Value2 = 1
$ws.Range("H124").Value2 = 1

# Row 130
$ws.Range("A130").Value2 = "Mauricio"
$ws.Range("C130").Value2 = 4

# Row 132
$ws.Range("A132").Value2 = "Etiopia"
$ws.Range("C132").Value2 = 1
